# Add 2022-Q3 data
#
# Semantics (derived from the target diff):
#   - Insert a new worksheet "2022-Q3" right after "总计" (i.e. at tab
#     position 2, pushing 2022-Q2/2022-Q1/2021-Q4/2021-Q2/2021-Q1/2020-Q4
#     one slot to the right - their content is untouched).
#   - Populate the new sheet with the 2022-Q3 per-fund holdings table
#     (same column layout / types as the other quarterly sheets).
#   - Insert a new summary row for "2022-Q3" right under the header row
#     of the "总计" sheet, shifting the existing quarters down and fixing
#     up the running index in column A.

$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item(1)
$insertBefore = $wb.Worksheets.Item(2)   # currently "2022-Q2" -> becomes position 3

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($insertBefore)
$q3.Name = "2022-Q3"

# Match the other quarterly sheets' page margins (0.75/0.75/1/1/0.5/0.5 in).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Header row (B1:H1) - same headers used by every quarterly sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Per-fund holdings for 2022-Q3 (index, code, name, size, position%, ratio%, value, rank).
# Columns C/D/E/F/G are stored as text in every quarterly sheet (not numbers),
# so we force the "@" text number format before assigning them.
$rows = @(
    @(0,  "519692", "交银成长混合A",                 "23.45", "76.71", "3.50", "0.8208", 8),
    @(1,  "519694", "交银蓝筹混合",                   "16.10", "78.09", "3.81", "0.6134", 8),
    @(2,  "000478", "建信中证500指数增强A",           "45.95", "82.53", "1.14", "0.5238", 2),
    @(3,  "460005", "华泰柏瑞价值增长混合A",          "9.00",  "93.35", "2.05", "0.1845", 7),
    @(4,  "008234", "光大保德信消费主题股票",         "2.63",  "90.48", "5.82", "0.1531", 6),
    @(5,  "005633", "建信中证500指数增强C",           "3.42",  "82.53", "1.14", "0.0390", 2),
    @(6,  "009726", "招商中证500等权重指数增强A",     "2.67",  "90.23", "1.44", "0.0384", 5),
    @(7,  "008778", "嘉实中证500指数增强A",           "1.08",  "91.33", "1.60", "0.0173", 10),
    @(8,  "009727", "招商中证500等权重指数增强C",     "1.12",  "90.23", "1.44", "0.0161", 5),
    @(9,  "008779", "嘉实中证500指数增强C",           "0.81",  "91.33", "1.60", "0.0130", 10),
    @(10, "003238", "新华外延增长主题灵活配置混合",   "0.50",  "57.43", "2.37", "0.0118", 5),
    @(11, "010154", "中加中证500指数增强C",           "0.51",  "94.15", "1.57", "0.0080", 7),
    @(12, "010153", "中加中证500指数增强A",           "0.44",  "94.15", "1.57", "0.0069", 7),
    @(13, "960016", "交银成长混合H",                 "0.16",  "76.71", "3.50", "0.0056", 8),
    @(14, "004546", "建信量化优享定期开放灵活配置混合", "0.16", "25.68", "1.06", "0.0017", 3),
    @(15, "010037", "华泰柏瑞价值增长混合C",          "0.05",  "93.35", "2.05", "0.0010", 7)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $row[0]          # A - running index (number)
    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[1]          # B - fund code (text)
    $q3.Cells.Item($r, 3).NumberFormat = "@"
    $q3.Cells.Item($r, 3).Value = $row[2]          # C - fund name (text)
    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[3]          # D - fund size (text)
    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[4]          # E - stock position (text)
    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[5]          # F - position ratio (text)
    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[6]          # G - held value (text)
    $q3.Cells.Item($r, 8).Value = $row[7]          # H - rank (number)
}

# Apply the same bold/bordered/centered style the other quarterly sheets
# use for their header row and index column (style index "2" in the
# original file - font w/o theme color, border w/o explicit color).
$zongji.Range("A2").Copy()
$q3.Range("A2:A17").PasteSpecial(-4122)

$zongji.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Insert the "2022-Q3" summary row into "总计", shifting the other
#    quarters down one row (and bumping their running index by 1).
# ---------------------------------------------------------------------
$oldRows = New-Object System.Collections.ArrayList
for ($r = 2; $r -le 7; $r++) {
    $b = $zongji.Cells.Item($r, 2).Value2
    $c = $zongji.Cells.Item($r, 3).Value2
    $d = $zongji.Cells.Item($r, 4).Value2
    [void]$oldRows.Add(@($b, $c, $d))
}

# New row 2: the 2022-Q3 summary.
$zongji.Cells.Item(2, 1).Value = 0
$zongji.Cells.Item(2, 2).Value = "2022-Q3"
$zongji.Cells.Item(2, 3).Value = 16
$zongji.Cells.Item(2, 4).Value = 2.45

# Push the previously-existing quarters (old rows 2-7) down into rows 3-8.
for ($i = 0; $i -lt $oldRows.Count; $i++) {
    $r = $i + 3
    $old = $oldRows[$i]
    $zongji.Cells.Item($r, 1).Value = $i + 1
    $zongji.Cells.Item($r, 2).Value = $old[0]
    $zongji.Cells.Item($r, 3).Value = $old[1]
    $zongji.Cells.Item($r, 4).Value = $old[2]
}

# Row 8 is brand new - give its index cell (A8) the same style as the
# rest of column A (style index "2").
$zongji.Range("A7").Copy()
$zongji.Range("A8").PasteSpecial(-4122)
$zongji.Cells.Item(8, 1).Value = 6

# ---------------------------------------------------------------------
# 3. Leave the view the same way it started: "总计" active, nothing else
#    marked as the selected tab.
# ---------------------------------------------------------------------
$zongji.Activate()
$zongji.Range("A1").Select()
